$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Formatting first (this never touches cell values / the shared-string
# table, so it can happen in any order): row 36's existing data switches
# from the "plain" style (4/5) to the "bottom border" style (8/9) because a
# new row (37) is appended right after it, closing off that translation
# block. Style 8/9 already exists in styles.xml (reused from rows 6/7/21),
# so copy *format only* from one of those rows instead of re-creating
# borders/fonts by hand - that keeps styles.xml untouched.
# ---------------------------------------------------------------------------
$ws.Range("A6:E6").Copy() | Out-Null
$ws.Range("A36:E36").PasteSpecial(-4122) | Out-Null

# New row 37 - single-line entry, same shape/style as rows 6/7/21 (8/9).
$ws.Range("A6:E6").Copy() | Out-Null
$ws.Range("A37:E37").PasteSpecial(-4122) | Out-Null
$ws.Rows("37:37").RowHeight = 43.2

# New row 38 - start of a new translation block (plain style 4/5, like row 2).
$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A38:E38").PasteSpecial(-4122) | Out-Null
$ws.Rows("38:38").RowHeight = 43.2

# New rows 39-41 - continuation rows, columns B-E only (plain style 4/5,
# like row 9).
$ws.Range("B9:E9").Copy() | Out-Null
$ws.Range("B39:E39").PasteSpecial(-4122) | Out-Null
$ws.Rows("39:39").RowHeight = 21.6

$ws.Range("B9:E9").Copy() | Out-Null
$ws.Range("B40:E40").PasteSpecial(-4122) | Out-Null
$ws.Rows("40:40").RowHeight = 21.6

$ws.Range("B9:E9").Copy() | Out-Null
$ws.Range("B41:E41").PasteSpecial(-4122) | Out-Null
$ws.Rows("41:41").RowHeight = 21.6

# ---------------------------------------------------------------------------
# Now fill in the values, in the same order the original author typed /
# pasted them (column by column), so brand-new shared-string entries land
# at the same indices as the target workbook.
# ---------------------------------------------------------------------------

# Row 37
$ws.Range("C37").Value = " I\'m glad to see your team back\nsafely, [hero]!"
$ws.Range("A37").Value = "SCRIPT/G01P04A/us2204.ssb"
$ws.Range("D37").Value = " Я рад, что твоя группа вернулась\nв целости и сохранности, [hero]!"
$ws.Range("E37").Value = " Ÿ ñàä, œóï óâïÿ ãñôððà âåñîôìàòû\nâ øåìïòóé é òïöñàîîïòóé, [hero]!"
$ws.Range("B37").Value = 49

# Rows 38-41, block A then column C, then column D, then column E.
$ws.Range("A38").Value = "SCRIPT/G01P04A/us2301.ssb"

$ws.Range("C38").Value = " Oh, [hero] and\n[partner]!"
$ws.Range("C39").Value = " I\'m happy to see you! Are you\nvisiting us?"
$ws.Range("C40").Value = " By the way...[K] Have you seen my\ndad anywhere?"
$ws.Range("C41").Value = " I haven\'t seen him for a while.\nI\'m beginning to get worried…"

$ws.Range("D38").Value = " О, [hero] и [partner]!"
$ws.Range("D39").Value = " Я рад вас видеть! Навещаете\nнас?"
$ws.Range("D40").Value = " Кстати...[K] Вы, случаем, не видели\nгде-нибудь моего отца?"
$ws.Range("D41").Value = " Я давно его не видел. Я начинаю\nбеспокоиться..."

$ws.Range("E38").Value = " Ï, [hero] é [partner]!"
$ws.Range("E39").Value = " Ÿ ñàä âàò âéäåóû! Îàâåþàåóå\nîàò?"
$ws.Range("E40").Value = " Ëòóàóé...[K] Âú, òìôœàåí, îå âéäåìé\nãäå-îéáôäû íïåãï ïóøà?"
$ws.Range("E41").Value = " Ÿ äàâîï åãï îå âéäåì. Ÿ îàœéîàý\náåòðïëïéóûòÿ..."

$ws.Range("B38").Value = 19
$ws.Range("B39").Value = 22
$ws.Range("B40").Value = 25
$ws.Range("B41").Value = 29

# ---------------------------------------------------------------------------
# Update the view: scroll so row 37 is at the top and select the last
# filled cell, matching where a user would land after typing this block.
# ---------------------------------------------------------------------------
$excel.Goto($ws.Range("A37"), $true) | Out-Null
$excel.ActiveWindow.ScrollRow = 37
$ws.Range("E41").Select() | Out-Null
